# "Zeitblätter" time-sheet update for Ursus Schneider — Oktober sheet.
# Adds two new daily entries (rows 20 & 21) describing DMX / OpenDMX / UDP
# testing work and documentation, then leaves the selection on D22 (as it
# was left in Excel after the entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")
$ws.Activate()

# Row 20 - Mittwoch (date already present via formula in A20)
$ws.Range("B20").Value = 4
$ws.Range("D20").Value = "14:00 - 15:30, 19:00 - 21:30"
$ws.Range("C20").Value = "DMX Tests in FH, OpenDMX test, UDP funktioniert"

# Row 21 - Donnerstag (date already present via formula in A21)
$ws.Range("B21").Value = 0.5
$ws.Range("C21").Value = "Dokumentation"
$ws.Range("D21").Value = "08:00 - 08:30"

# Leave the cursor where the author left it after typing the new rows.
$ws.Range("D22").Select()
